# Update column G ("K" strikeout count) values for rows 2-23 on the active sheet.
# Mirrors: regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 3
    3  = 3
    4  = 3
    5  = 1
    6  = 1
    7  = 0
    8  = 1
    9  = 1
    10 = 0
    11 = 0
    12 = 1
    13 = 1
    14 = 1
    15 = 0
    16 = 1
    17 = 3
    18 = 1
    19 = 1
    20 = 0
    21 = 7
    22 = 2
    23 = 4
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
